$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Add the new "Purchase Order Line Items" sheet
#    (moved to the end of the workbook only after all writes below are
#    done - re-indexing on Move invalidates the $ws handle for further use)
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Add()
$ws.Name = "Purchase Order Line Items"

# Column widths (A:8, B:12, C:10, D:11, E:19) - ColumnWidth needs the ~0.8333
# padding subtracted so the stored sheet width lands on the exact integer.
$pad = 0.8333333333333
$ws.Columns.Item(1).ColumnWidth = 8 - $pad
$ws.Columns.Item(2).ColumnWidth = 12 - $pad
$ws.Columns.Item(3).ColumnWidth = 10 - $pad
$ws.Columns.Item(4).ColumnWidth = 11 - $pad
$ws.Columns.Item(5).ColumnWidth = 19 - $pad

# Header row
$headers = @("PO ID", "Product ID", "Quantity", "Unit Cost", "Quantity Received")
for ($col = 1; $col -le $headers.Count; $col++) {
    $ws.Cells.Item(1, $col).Value = $headers[$col - 1]
}

# Bold, centered, light-blue header styling (matches Expenses/Revenue modals)
$headerRange = $ws.Range("A1:E1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$lightBlue = 0xAD + (0xD8 * 256) + (0xE6 * 65536)
$headerRange.Interior.Color = $lightBlue

# Line item rows: PO ID, Product ID, Quantity, Unit Cost, Quantity Received
$rows = @(
    @("PO-001", "PRD-001", 5,  899,  5),
    @("PO-001", "PRD-004", 5,  699,  5),
    @("PO-001", "PRD-007", 11, 149,  11),
    @("PO-002", "PRD-002", 3,  1349, 3),
    @("PO-002", "PRD-006", 6,  349,  6),
    @("PO-003", "PRD-008", 2,  1195, 2),
    @("PO-003", "PRD-009", 1,  1515, 1),
    @("PO-004", "PRD-002", 4,  1349, 4),
    @("PO-004", "PRD-001", 3,  899,  3),
    @("PO-004", "PRD-007", 4,  149,  4),
    @("PO-005", "PRD-005", 3,  1999, 3),
    @("PO-005", "PRD-004", 2,  699,  2),
    @("PO-005", "PRD-007", 3,  149,  3),
    @("PO-006", "PRD-001", 4,  899,  4),
    @("PO-006", "PRD-006", 5,  349,  5),
    @("PO-006", "PRD-007", 3,  149,  3),
    @("PO-007", "PRD-002", 3,  1349, 0),
    @("PO-007", "PRD-004", 2,  699,  0),
    @("PO-007", "PRD-006", 3,  349,  0),
    @("PO-008", "PRD-003", 6,  849,  0)
)

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $rows[$i]
    $rowNum = $i + 2
    for ($col = 1; $col -le $r.Count; $col++) {
        $ws.Cells.Item($rowNum, $col).Value = $r[$col - 1]
    }
}

# Now move the finished sheet to the end of the tab order.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws.Move($null, $lastSheet)

# ---------------------------------------------------------------------------
# 2. Clean up stray empty cells in the Products sheet (J17:J21)
# ---------------------------------------------------------------------------
$products = $wb.Worksheets.Item("Products")
$products.Range("J17:J21").ClearContents()
